# Add 2022-Q3 data
# 1. Update the "总计" (Total) summary sheet: insert a new leading data row for
#    2022-Q3 and push the existing 2021-Q2 / 2020-Q4 rows down by one.
# 2. Insert a brand-new worksheet named "2022-Q3" right after the Total sheet,
#    holding the per-fund holding detail for that quarter.

$wb = $excel.ActiveWorkbook
$total = $wb.Worksheets.Item(1)

# --- Step 1: shift rows 2-3 of the Total sheet down to 3-4, then write the
#     new 2022-Q3 values into row 2. Values are written cell-by-cell (instead
#     of a row Insert) so no incidental row-level style bleeds onto the new
#     cells; the single style-carrying cell (column A) is copied explicitly.

# Preserve column-A styling for the row that is about to become row 4 by
# copying it from the existing row 3 (style only matters, value is
# overwritten immediately after).
$total.Cells.Item(3, 1).Copy($total.Cells.Item(4, 1))

# Old row 3 (2020-Q4) -> row 4
$total.Cells.Item(4, 1).Value = 2
$total.Cells.Item(4, 2).Value = "2020-Q4"
$total.Cells.Item(4, 3).Value = 2
$total.Cells.Item(4, 4).Value = 0.02

# Old row 2 (2021-Q2) -> row 3
$total.Cells.Item(3, 1).Value = 1
$total.Cells.Item(3, 2).Value = "2021-Q2"
$total.Cells.Item(3, 3).Value = 2
$total.Cells.Item(3, 4).Value = 0.01

# New row 2: 2022-Q3
$total.Cells.Item(2, 1).Value = 0
$total.Cells.Item(2, 2).Value = "2022-Q3"
$total.Cells.Item(2, 3).Value = 3
$total.Cells.Item(2, 4).Value = 0.11

# --- Step 2: insert the new "2022-Q3" worksheet right after the Total sheet.
$q3 = $wb.Worksheets.Add($null, $total)
$q3.Name = "2022-Q3"

# Columns B (基金代码) and D:G (规模/仓位/市值 figures) hold digit-only text
# such as "002291" or "0.0928" in the source data. A plain .Value assignment
# would auto-coerce those numeric-looking strings to real numbers (dropping
# leading zeros / introducing float noise), so mark the ranges as Text first
# to force a faithful string write, matching the original inlineStr cells.
$q3.Range("B2:B4").NumberFormat = "@"
$q3.Range("D2:G4").NumberFormat = "@"

# Header row
$q3.Cells.Item(1, 2).Value = "基金代码"
$q3.Cells.Item(1, 3).Value = "基金名称"
$q3.Cells.Item(1, 4).Value = "基金规模"
$q3.Cells.Item(1, 5).Value = "股票总仓位"
$q3.Cells.Item(1, 6).Value = "仓位占比"
$q3.Cells.Item(1, 7).Value = "持有市值(亿元)"
$q3.Cells.Item(1, 8).Value = "仓位排名"

# Row 2
$q3.Cells.Item(2, 1).Value = 0
$q3.Cells.Item(2, 2).Value = "002291"
$q3.Cells.Item(2, 3).Value = "诺安安鑫灵活配置混合"
$q3.Cells.Item(2, 4).Value = "2.66"
$q3.Cells.Item(2, 5).Value = "77.38"
$q3.Cells.Item(2, 6).Value = "3.49"
$q3.Cells.Item(2, 7).Value = "0.0928"
$q3.Cells.Item(2, 8).Value = 6

# Row 3
$q3.Cells.Item(3, 1).Value = 1
$q3.Cells.Item(3, 2).Value = "002137"
$q3.Cells.Item(3, 3).Value = "诺安利鑫灵活配置混合A"
$q3.Cells.Item(3, 4).Value = "0.44"
$q3.Cells.Item(3, 5).Value = "76.46"
$q3.Cells.Item(3, 6).Value = "3.61"
$q3.Cells.Item(3, 7).Value = "0.0159"
$q3.Cells.Item(3, 8).Value = 5

# Row 4
$q3.Cells.Item(4, 1).Value = 2
$q3.Cells.Item(4, 2).Value = "014521"
$q3.Cells.Item(4, 3).Value = "诺安利鑫灵活配置混合C"
$q3.Cells.Item(4, 4).Value = "0.05"
$q3.Cells.Item(4, 5).Value = "76.46"
$q3.Cells.Item(4, 6).Value = "3.61"
$q3.Cells.Item(4, 7).Value = "0.0018"
$q3.Cells.Item(4, 8).Value = 5
